# Perturbation: reshape the "optimization_parameters" sheet (new
# "production_function" / "L_curve" rows replacing the old "Model" label
# and the stray "Deletion" row) and switch the active sheet from
# "optimization_diagnostics" to "optimization_parameters".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# --- Row 1: drop the extra duplicated "value" header cells (C1:F1) -----
$ws.Range("C1:F1").ClearContents()

# --- Make room for the new "L_curve" row right after the model row -----
# (old layout: row8 = Model/Sigmoid, row9 = estimate_params, ...
#  row16 = Deletion/0/3, row17 = simulation_timepoints)
$ws.Rows.Item(9).Insert()

# Row 8 label changes from "Model" to "production_function" (value stays
# "Sigmoid").
$ws.Range("A8").Value = "production_function"

# New row 9: "L_curve" flag, defaulting to 1, formatted like the other
# numeric parameter values (scientific notation, same as B2/B4/B5/B6/B7).
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 1
$ws.Range("B9").NumberFormat = "0.00E+00"

# --- Remove the obsolete "Deletion" row (now shifted down to row 17, ---
# --- just above the "simulation_timepoints" row) ------------------------
$ws.Rows.Item(17).Delete()

# --- Switch the active sheet + selection to optimization_parameters ----
$ws.Activate()
$ws.Rows.Item(17).Select()
